$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as TEXT, preserving the default (no explicit) style
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '67.476.87'
Set-TextCell $ws.Range("E2") '  -0.33%  '

# Row 3
Set-TextCell $ws.Range("D3") '3.778.74'
Set-TextCell $ws.Range("E3") '  -0.40%  '

# Row 4
Set-TextCell $ws.Range("D4") '1.00'
Set-TextCell $ws.Range("E4") '  +0.00%  '

# Row 5
Set-TextCell $ws.Range("D5") '597.34'
Set-TextCell $ws.Range("E5") '  +0.27%  '

# Row 6
Set-TextCell $ws.Range("D6") '164.36'
Set-TextCell $ws.Range("E6") '  -1.60%  '

# Row 8
Set-TextCell $ws.Range("E8") '  -1.05%  '

# Row 9
Set-TextCell $ws.Range("E9") '  -1.13%  '

# Row 10
Set-TextCell $ws.Range("D10") '0.448'
Set-TextCell $ws.Range("E10") '  +0.11%  '

# Row 11
Set-TextCell $ws.Range("E11") '  +1.44%  '

# Row 12
Set-TextCell $ws.Range("E12") '  -2.32%  '

# Row 13
Set-TextCell $ws.Range("D13") '35.54'
Set-TextCell $ws.Range("E13") '  -1.51%  '

# Row 14
Set-TextCell $ws.Range("D14") '4.416.16'
Set-TextCell $ws.Range("E14") '  -0.29%  '

# Row 15
Set-TextCell $ws.Range("D15") '3.758.30'
Set-TextCell $ws.Range("E15") '  -1.04%  '

# Row 16
Set-TextCell $ws.Range("D16") '67.563.41'
Set-TextCell $ws.Range("E16") '  -0.15%  '

# Row 17
Set-TextCell $ws.Range("D17") '18.23'
Set-TextCell $ws.Range("E17") '  -0.91%  '

# Row 18
Set-TextCell $ws.Range("E18") '  +1.68%  '

# Row 19
Set-TextCell $ws.Range("D19") '7.00'
Set-TextCell $ws.Range("E19") '  -0.83%  '

# Row 20
Set-TextCell $ws.Range("D20") '460.13'
Set-TextCell $ws.Range("E20") '  +0.19%  '

# Row 21
Set-TextCell $ws.Range("E21") '  -2.40%  '

# Row 22
Set-TextCell $ws.Range("D22") '0.694'
Set-TextCell $ws.Range("E22") '  -0.32%  '

# Row 23
Set-TextCell $ws.Range("E23") '  -6.07%  '

# Row 24
Set-TextCell $ws.Range("D24") '82.32'
Set-TextCell $ws.Range("E24") '  -1.19%  '

# Row 25
Set-TextCell $ws.Range("D25") '11.95'
Set-TextCell $ws.Range("E25") '  -0.86%  '

# Row 26
Set-TextCell $ws.Range("E26") '  -1.53%  '

# Row 27
Set-TextCell $ws.Range("E27") '  -0.03%  '

# Row 28
Set-TextCell $ws.Range("D28") '9.96'
Set-TextCell $ws.Range("E28") '  -0.50%  '

# Row 29
Set-TextCell $ws.Range("D29") '3.927.97'
Set-TextCell $ws.Range("E29") '  -0.34%  '

# Row 30
Set-TextCell $ws.Range("D30") '7.39'
Set-TextCell $ws.Range("E30") '  +2.01%  '

# Row 31
Set-TextCell $ws.Range("D31") '2.64'
Set-TextCell $ws.Range("E31") '  -4.57%  '

# Row 32
Set-TextCell $ws.Range("E32") '  -2.74%  '

# Row 33
Set-TextCell $ws.Range("D33") '28.89'
Set-TextCell $ws.Range("E33") '  -2.59%  '

# Row 34
Set-TextCell $ws.Range("D34") '0.996'
Set-TextCell $ws.Range("E34") '  -0.35%  '

# Row 35
Set-TextCell $ws.Range("D35") '8.95'

# Row 36
Set-TextCell $ws.Range("D36") '0.0984'
Set-TextCell $ws.Range("E36") '  -1.70%  '

# Row 37
Set-TextCell $ws.Range("E37") '  +0.03%  '

# Row 38
Set-TextCell $ws.Range("D38") '0.988'
Set-TextCell $ws.Range("E38") '  -0.29%  '

# Row 39
Set-TextCell $ws.Range("D39") '3.22'
Set-TextCell $ws.Range("E39") '  -4.76%  '

# Row 40
Set-TextCell $ws.Range("E40") '  -0.79%  '

# Row 41
Set-TextCell $ws.Range("E41") '  +0.09%  '

# Row 43
Set-TextCell $ws.Range("D43") '43.58'
Set-TextCell $ws.Range("E43") '  -1.43%  '

# Row 44
Set-TextCell $ws.Range("D44") '47.47'
Set-TextCell $ws.Range("E44") '  -1.21%  '

# Row 45
Set-TextCell $ws.Range("D45") '0.295'
Set-TextCell $ws.Range("E45") '  -0.85%  '

# Row 46
Set-TextCell $ws.Range("D46") '150.90'
Set-TextCell $ws.Range("E46") '  +0.92%  '

# Row 47
Set-TextCell $ws.Range("E47") '  +0.25%  '

# Row 48
Set-TextCell $ws.Range("E48") '  +7.61%  '

# Row 49
Set-TextCell $ws.Range("D49") '27.03'
Set-TextCell $ws.Range("E49") '  +0.95%  '

# Rows 50 and 51: swap Bittensor and Stacks entries with updated values
Set-TextCell $ws.Range("B50") 'Stacks'
Set-TextCell $ws.Range("C50") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws.Range("D50") '1.84'
Set-TextCell $ws.Range("E50") '  +1.12%  '

Set-TextCell $ws.Range("B51") 'Bittensor'
Set-TextCell $ws.Range("C51") 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws.Range("D51") '391.59'
Set-TextCell $ws.Range("E51") '  -0.44%  '
